$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings are not
# auto-converted to numbers by Excel (matches source data which stores
# these as literal text).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.025.59"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.642.49"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "214.72"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "0.5096"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.2566"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "0.06356"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "19.56"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "0.07757"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "4.282"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "1.640.00"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "64.33"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "0.0₅7718"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "26.053.90"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "197.16"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "4.419"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "9.924"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "6.027"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "1.861"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "141.46"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").Value = "0.1191"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("D27").Value = "6.815"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "15.59"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "0.04851"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "3.248"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "3.163"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "1.524"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "0.8961"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "2.580"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "1.139.64"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "0.5446"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "2.521"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("D43").Value = "0.8093"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "99.34"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "5.393"
$ws.Range("E45").Value = "  -5.25%  "
$ws.Range("D46").Value = "1.780.24"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "0.4529"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "0.9989"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "54.80"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "0.05059"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("E51").Value = "  -0.48%  "

# Restore the default (Normal) cell style now that the text values are
# set, so no lingering custom number-format style is left on the cells.
$dRange.Style = "Normal"
